$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: runs=6, balls=4, fours=1
$ws.Range("C2").Value = "'6"
$ws.Range("D2").Value = "'4"
$ws.Range("E2").Value = "'1"

# Row 4: runs=2, balls=3, fours=0
$ws.Range("C4").Value = "'2"
$ws.Range("D4").Value = "'3"
$ws.Range("E4").Value = "'0"

# Row 5: runs=1, balls=2
$ws.Range("C5").Value = "'1"
$ws.Range("D5").Value = "'2"

# Row 6: runs=5, balls=7
$ws.Range("C6").Value = "'5"
$ws.Range("D6").Value = "'7"
